# Applies a weekly re-shuffle of the per-row price/date/variety data for the
# "Hortaliza, Terminal La Palmera de La Serena - Alcachofa" sheet.
#
# The edit moves the (Fecha, Variedad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades) tuple that used to live on one data row onto another data
# row, i.e. it is a permutation of the existing rows 2..26 restricted to the
# columns D, H, J, K, L, M, N, O, P, Q (every other column - Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Calidad, Clasificacion -
# stays constant for every row so it does not need to move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry row-specific data that gets permuted.
$cols = @("D", "H", "J", "K", "L", "M", "N", "O", "P", "Q")

# new row number -> old row number the data comes from.
$mapping = @{
    2  = 19
    3  = 5
    4  = 6
    5  = 21
    6  = 23
    7  = 3
    8  = 22
    9  = 14
    10 = 10
    11 = 11
    12 = 24
    13 = 12
    14 = 13
    15 = 16
    16 = 20
    17 = 7
    18 = 8
    19 = 15
    20 = 4
    21 = 18
    22 = 25
    23 = 17
    24 = 26
    25 = 2
    26 = 9
}

# Snapshot every source cell's value before any writes happen, since some
# rows both receive and donate data (the permutation has cycles longer than
# a simple swap), so writing in row order would clobber data we still need.
# `.Value2` (not `.Value`) is used because it is the property that actually
# round-trips scalars (numbers/dates/strings) through this COM bridge.
$snapshot = @{}
foreach ($row in 2..26) {
    foreach ($col in $cols) {
        $snapshot["$col$row"] = $ws.Range("$col$row").Value2
    }
}

foreach ($row in 2..26) {
    $srcRow = $mapping[$row]
    if ($srcRow -eq $row) {
        continue
    }
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $snapshot["$col$srcRow"]
    }
}
